$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1:G21").ClearContents()

# --- Header row (order: A1,B1,C1,D1,E1,F1,G1) matches original shared string order 0-5 ---
$ws.Cells.Item(1, 1).Value2 = "document"
$ws.Cells.Item(1, 2).Value2 = "url"
$ws.Cells.Item(1, 3).Value2 = "h1"
$ws.Cells.Item(1, 4).Value2 = "published"
$ws.Cells.Item(1, 5).Value2 = "comment"
$ws.Cells.Item(1, 6).Value2 = "document"
$ws.Cells.Item(1, 7).Value2 = "updated"

# --- Column A (rows 2-21) ---
$ws.Cells.Item(2, 1).Value2 = "What is the Technology Window"
$ws.Cells.Item(3, 1).Value2 = "What is the Skyscraper Method"
$ws.Cells.Item(4, 1).Value2 = "What is the single source of truth principle"
$ws.Cells.Item(5, 1).Value2 = "What is the self-fulfilling prophecy"
$ws.Cells.Item(6, 1).Value2 = "What is the quantitative accumulation effect"
$ws.Cells.Item(7, 1).Value2 = "What is the Narrow Path Principle"
$ws.Cells.Item(8, 1).Value2 = "What is the Dunning-Krueger-Effect"
$ws.Cells.Item(9, 1).Value2 = "What is the Cheerleader effect"
$ws.Cells.Item(10, 1).Value2 = "What is the Bystander Effect"
$ws.Cells.Item(11, 1).Value2 = "What is the believe bias"
$ws.Cells.Item(12, 1).Value2 = "What is the barnum effect"
$ws.Cells.Item(13, 1).Value2 = "What is the ambiguity effect"
$ws.Cells.Item(14, 1).Value2 = "What is the 80 percent rule"
$ws.Cells.Item(15, 1).Value2 = "What is an Aha Moment"
$ws.Cells.Item(16, 1).Value2 = "13 cognitive biases to avoid in daily life"
$ws.Cells.Item(17, 1).Value2 = "Key setup for problem solving"
$ws.Cells.Item(18, 1).Value2 = "What is K Level Thinking"
$ws.Cells.Item(19, 1).Value2 = "How can devops team take advantage of artificial intelligence "
$ws.Cells.Item(20, 1).Value2 = "What Jobs are threatened by AI the most"
$ws.Cells.Item(21, 1).Value2 = "How to Break Down a Complex Problem"

# --- Column B (rows 2-21) ---
$ws.Cells.Item(2, 2).Value2 = "what-is-the-technology-window"
$ws.Cells.Item(3, 2).Value2 = "what-is-the-skyscraper-method"
$ws.Cells.Item(4, 2).Value2 = "what-is-the-single-source-of-truth-principle"
$ws.Cells.Item(5, 2).Value2 = "what-is-the-self-fulfilling-prophecy"
$ws.Cells.Item(6, 2).Value2 = "what-is-the-quantitative-accumulation-effect"
$ws.Cells.Item(7, 2).Value2 = "what-is-the-narrow-path-principle"
$ws.Cells.Item(8, 2).Value2 = "what-is-the-dunning-krueger-effect"
$ws.Cells.Item(9, 2).Value2 = "what-is-the-cheerleader-effect"
$ws.Cells.Item(10, 2).Value2 = "what-is-the-bystander-effect"
$ws.Cells.Item(11, 2).Value2 = "what-is-the-believe-bias"
$ws.Cells.Item(12, 2).Value2 = "what-is-the-barnum-effect"
$ws.Cells.Item(13, 2).Value2 = "what-is-the-ambiguity-effect"
$ws.Cells.Item(14, 2).Value2 = "what-is-the-80-percent-rule"
$ws.Cells.Item(15, 2).Value2 = "what-is-an-aha-moment"
$ws.Cells.Item(16, 2).Value2 = "cognitive-biases-to-avoid-in-daily-life"
$ws.Cells.Item(17, 2).Value2 = "key-setup-for-problem-solving"
$ws.Cells.Item(18, 2).Value2 = "what-is-k-level-thinking"
$ws.Cells.Item(19, 2).Value2 = "how-can-devops-team-take-advantage-of-artificial-intelligence-"
$ws.Cells.Item(20, 2).Value2 = "what-jobs-are-threatened-by-ai-the-most"
$ws.Cells.Item(21, 2).Value2 = "how-to-break-down-a-complex-problem"

# --- Column C (rows 2-21) ---
$ws.Cells.Item(2, 3).Value2 = "What Is The Technology Window"
$ws.Cells.Item(3, 3).Value2 = "What Is The Skyscraper Method"
$ws.Cells.Item(4, 3).Value2 = "What Is The Single Source Of Truth Principle"
$ws.Cells.Item(5, 3).Value2 = "What Is The Self-Fulfilling Prophecy"
$ws.Cells.Item(6, 3).Value2 = "What Is The Quantitative Accumulation Effect"
$ws.Cells.Item(7, 3).Value2 = "What Is The Narrow Path Principle"
$ws.Cells.Item(8, 3).Value2 = "What Is The Dunning-Krueger-Effect"
$ws.Cells.Item(9, 3).Value2 = "What Is The Cheerleader Effect"
$ws.Cells.Item(10, 3).Value2 = "What Is The Bystander Effect"
$ws.Cells.Item(11, 3).Value2 = "What Is The Believe Bias"
$ws.Cells.Item(12, 3).Value2 = "What Is The Barnum Effect"
$ws.Cells.Item(13, 3).Value2 = "What Is The Ambiguity Effect"
$ws.Cells.Item(14, 3).Value2 = "What Is The 80 Percent Rule"
$ws.Cells.Item(15, 3).Value2 = "What Is An Aha Moment"
$ws.Cells.Item(16, 3).Value2 = "13 Cognitive Biases To Avoid In Daily Life"
$ws.Cells.Item(17, 3).Value2 = "Key Setup For Problem Solving"
$ws.Cells.Item(18, 3).Value2 = "What Is K Level Thinking"
$ws.Cells.Item(19, 3).Value2 = "How Can Devops Team Take Advantage Of Artificial Intelligence "
$ws.Cells.Item(20, 3).Value2 = "What Jobs Are Threatened By Ai The Most"
$ws.Cells.Item(21, 3).Value2 = "How To Break Down A Complex Problem"

# --- Column F (rows 2-21), duplicates column A strings ---
$ws.Cells.Item(2, 6).Value2 = "What is the Technology Window"
$ws.Cells.Item(3, 6).Value2 = "What is the Skyscraper Method"
$ws.Cells.Item(4, 6).Value2 = "What is the single source of truth principle"
$ws.Cells.Item(5, 6).Value2 = "What is the self-fulfilling prophecy"
$ws.Cells.Item(6, 6).Value2 = "What is the quantitative accumulation effect"
$ws.Cells.Item(7, 6).Value2 = "What is the Narrow Path Principle"
$ws.Cells.Item(8, 6).Value2 = "What is the Dunning-Krueger-Effect"
$ws.Cells.Item(9, 6).Value2 = "What is the Cheerleader effect"
$ws.Cells.Item(10, 6).Value2 = "What is the Bystander Effect"
$ws.Cells.Item(11, 6).Value2 = "What is the believe bias"
$ws.Cells.Item(12, 6).Value2 = "What is the barnum effect"
$ws.Cells.Item(13, 6).Value2 = "What is the ambiguity effect"
$ws.Cells.Item(14, 6).Value2 = "What is the 80 percent rule"
$ws.Cells.Item(15, 6).Value2 = "What is an Aha Moment"
$ws.Cells.Item(16, 6).Value2 = "13 cognitive biases to avoid in daily life"
$ws.Cells.Item(17, 6).Value2 = "Key setup for problem solving"
$ws.Cells.Item(18, 6).Value2 = "What is K Level Thinking"
$ws.Cells.Item(19, 6).Value2 = "How can devops team take advantage of artificial intelligence "
$ws.Cells.Item(20, 6).Value2 = "What Jobs are threatened by AI the most"
$ws.Cells.Item(21, 6).Value2 = "How to Break Down a Complex Problem"

# --- Column D/G numeric (rows 2-21) ---
$ws.Cells.Item(2, 4).Value2 = 44984.79226117156
$ws.Cells.Item(2, 7).Value2 = 44968.82226430001
$ws.Cells.Item(3, 4).Value2 = 44984.79139847053
$ws.Cells.Item(3, 7).Value2 = 44968.82335802943
$ws.Cells.Item(4, 4).Value2 = 44984.79075432817
$ws.Cells.Item(4, 7).Value2 = 44968.8283912778
$ws.Cells.Item(5, 4).Value2 = 44984.78997172293
$ws.Cells.Item(5, 7).Value2 = 44968.83482908575
$ws.Cells.Item(6, 4).Value2 = 44984.78734691756
$ws.Cells.Item(6, 7).Value2 = 44968.82762568397
$ws.Cells.Item(7, 4).Value2 = 44984.78653146478
$ws.Cells.Item(7, 7).Value2 = 44968.836970087
$ws.Cells.Item(8, 4).Value2 = 44984.77886384894
$ws.Cells.Item(8, 7).Value2 = 44968.83607997302
$ws.Cells.Item(9, 4).Value2 = 44984.76482031686
$ws.Cells.Item(9, 7).Value2 = 44968.83272623207
$ws.Cells.Item(10, 4).Value2 = 44984.76419652337
$ws.Cells.Item(10, 7).Value2 = 44968.83799102488
$ws.Cells.Item(11, 4).Value2 = 44984.76331853376
$ws.Cells.Item(11, 7).Value2 = 44968.83201703853
$ws.Cells.Item(12, 4).Value2 = 44984.76202055444
$ws.Cells.Item(12, 7).Value2 = 44968.83331997384
$ws.Cells.Item(13, 4).Value2 = 44984.76082692308
$ws.Cells.Item(13, 7).Value2 = 44968.83411760339
$ws.Cells.Item(14, 4).Value2 = 44984.75974282543
$ws.Cells.Item(14, 7).Value2 = 44968.82652412132
$ws.Cells.Item(15, 4).Value2 = 44983.66443622751
$ws.Cells.Item(15, 7).Value2 = 44968.82335802943
$ws.Cells.Item(16, 4).Value2 = 44983.66333873353
$ws.Cells.Item(16, 7).Value2 = 44969.74406943242
$ws.Cells.Item(17, 4).Value2 = 44970.84635276953
$ws.Cells.Item(17, 7).Value2 = 44969.9128825571
$ws.Cells.Item(18, 4).Value2 = 44969.5374997161
$ws.Cells.Item(18, 7).Value2 = 44968.81650771474
$ws.Cells.Item(19, 4).Value2 = 44969.4746938719
$ws.Cells.Item(19, 7).Value2 = 44969.40788246163
$ws.Cells.Item(20, 4).Value2 = 44968.80385106849
$ws.Cells.Item(20, 7).Value2 = 44968.42467377915
$ws.Cells.Item(21, 4).Value2 = 44963.57076442974
$ws.Cells.Item(21, 7).Value2 = 44962.65609726021
